# Applies the "scheduled runner" value updates to the Phantom_Profits workbook.
# Each sheet corresponds to a crafting job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR);
# the diff only touches pre-existing numeric columns H-N (price/profit columns).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3134.4138
$ws.Range("I40").Value = 1694.3889
$ws.Range("K40").Value = 1694.3889
$ws.Range("M40").Value = -1519.3889
$ws.Range("H69").Value = 14849.875
$ws.Range("I69").Value = 14759.8
$ws.Range("K69").Value = 44279.39999999999
$ws.Range("M69").Value = -43405.39999999999
$ws.Range("H72").Value = 14849.875
$ws.Range("I72").Value = 14759.8
$ws.Range("K72").Value = 132838.2
$ws.Range("M72").Value = -128470.2
$ws.Range("H95").Value = 100624
$ws.Range("J95").Value = 100624
$ws.Range("L95").Value = 100624
$ws.Range("N95").Value = -106116
$ws.Range("H98").Value = 1744.95
$ws.Range("I98").Value = 1649.1428
$ws.Range("J98").Value = 1968.5
$ws.Range("K98").Value = 1649.1428
$ws.Range("L98").Value = 1968.5
$ws.Range("M98").Value = -151.1428000000001
$ws.Range("N98").Value = -4964.5
$ws.Range("H116").Value = 5799.2856
$ws.Range("I116").Value = 5266
$ws.Range("K116").Value = 5266
$ws.Range("M116").Value = -1824
$ws.Range("H118").Value = 2044.5
$ws.Range("I118").Value = 2044.5
$ws.Range("K118").Value = 6133.5
$ws.Range("M118").Value = -4476.5
$ws.Range("H122").Value = 1744.95
$ws.Range("I122").Value = 1649.1428
$ws.Range("J122").Value = 1968.5
$ws.Range("K122").Value = 4947.428400000001
$ws.Range("L122").Value = 5905.5
$ws.Range("M122").Value = -2497.428400000001
$ws.Range("N122").Value = -10805.5
$ws.Range("H132").Value = 3001.5186
$ws.Range("I132").Value = 2578.3809
$ws.Range("K132").Value = 7735.1427
$ws.Range("M132").Value = -5205.1427
$ws.Range("H135").Value = 548
$ws.Range("I135").Value = 420.33334
$ws.Range("J135").Value = 931
$ws.Range("K135").Value = 3783.00006
$ws.Range("L135").Value = 8379
$ws.Range("M135").Value = -1248.00006
$ws.Range("N135").Value = -13449

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10000
$ws.Range("I61").Value = 10000
$ws.Range("K61").Value = 10000
$ws.Range("M61").Value = -9788
$ws.Range("H74").Value = 3426.0527
$ws.Range("I74").Value = 2521.4167
$ws.Range("K74").Value = 2521.4167
$ws.Range("M74").Value = -1647.4167
$ws.Range("H77").Value = 3426.0527
$ws.Range("I77").Value = 2521.4167
$ws.Range("K77").Value = 12607.0835
$ws.Range("M77").Value = -8239.083500000001
$ws.Range("H97").Value = 2285.4285
$ws.Range("J97").Value = 3935
$ws.Range("L97").Value = 3935
$ws.Range("N97").Value = -4927
$ws.Range("H125").Value = 72999.664
$ws.Range("J125").Value = 72999.664
$ws.Range("L125").Value = 72999.664
$ws.Range("N125").Value = -82839.664
$ws.Range("H132").Value = 4200.9
$ws.Range("I132").Value = 3549.7
$ws.Range("K132").Value = 10649.1
$ws.Range("M132").Value = -8119.099999999999
$ws.Range("H136").Value = 10000
$ws.Range("I136").Value = 10000
$ws.Range("K136").Value = 30000
$ws.Range("M136").Value = -27450

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4060.6
$ws.Range("I20").Value = 4060.6
$ws.Range("K20").Value = 4060.6
$ws.Range("M20").Value = -3813.6
$ws.Range("H86").Value = 11760.75
$ws.Range("I86").Value = 7077.25
$ws.Range("J86").Value = 16444.25
$ws.Range("K86").Value = 7077.25
$ws.Range("L86").Value = 16444.25
$ws.Range("M86").Value = -5954.25
$ws.Range("N86").Value = -18690.25
$ws.Range("H89").Value = 11760.75
$ws.Range("I89").Value = 7077.25
$ws.Range("J89").Value = 16444.25
$ws.Range("K89").Value = 35386.25
$ws.Range("L89").Value = 82221.25
$ws.Range("M89").Value = -29770.25
$ws.Range("N89").Value = -93453.25
$ws.Range("H96").Value = 24948.5
$ws.Range("I96").Value = 24948.5
$ws.Range("K96").Value = 24948.5
$ws.Range("M96").Value = -22202.5
$ws.Range("H99").Value = 1452.7142
$ws.Range("I99").Value = 1045
$ws.Range("K99").Value = 1045
$ws.Range("M99").Value = 453
$ws.Range("H134").Value = 400
$ws.Range("I134").Value = 400
$ws.Range("K134").Value = 1200
$ws.Range("M134").Value = 1335

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1122.7142
$ws.Range("I16").Value = 1122.7142
$ws.Range("K16").Value = 1122.7142
$ws.Range("M16").Value = -835.7141999999999
$ws.Range("H41").Value = 9000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H58").Value = 2049
$ws.Range("I58").Value = 1452.4
$ws.Range("K58").Value = 1452.4
$ws.Range("M58").Value = -1249.4
$ws.Range("H99").Value = 5973.143
$ws.Range("I99").Value = 5973.143
$ws.Range("K99").Value = 5973.143
$ws.Range("M99").Value = -4475.143
$ws.Range("H107").Value = 633.619
$ws.Range("I107").Value = 521.4211
$ws.Range("K107").Value = 521.4211
$ws.Range("M107").Value = 1398.5789
$ws.Range("H113").Value = 1122.7142
$ws.Range("I113").Value = 1122.7142
$ws.Range("K113").Value = 1122.7142
$ws.Range("M113").Value = 1047.2858
$ws.Range("H126").Value = 5973.143
$ws.Range("I126").Value = 5973.143
$ws.Range("K126").Value = 17919.429
$ws.Range("M126").Value = -15449.429
$ws.Range("H134").Value = 1835.6666
$ws.Range("I134").Value = 1516.2667
$ws.Range("K134").Value = 4548.800099999999
$ws.Range("M134").Value = -2013.800099999999
$ws.Range("H136").Value = 2049
$ws.Range("I136").Value = 1452.4
$ws.Range("K136").Value = 4357.200000000001
$ws.Range("M136").Value = -1807.200000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 837.6667
$ws.Range("I68").Value = 733.3333
$ws.Range("J68").Value = 942
$ws.Range("K68").Value = 2199.9999
$ws.Range("L68").Value = 2826
$ws.Range("M68").Value = -1388.9999
$ws.Range("N68").Value = -4448
$ws.Range("H71").Value = 837.6667
$ws.Range("I71").Value = 733.3333
$ws.Range("J71").Value = 942
$ws.Range("K71").Value = 6599.9997
$ws.Range("L71").Value = 8478
$ws.Range("M71").Value = -2543.9997
$ws.Range("N71").Value = -16590
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 2624.75
$ws.Range("I131").Value = 1999.8
$ws.Range("K131").Value = 5999.4
$ws.Range("M131").Value = -959.3999999999996

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 25000
$ws.Range("I28").Value = 15000
$ws.Range("J28").Value = 35000
$ws.Range("K28").Value = 15000
$ws.Range("L28").Value = 35000
$ws.Range("M28").Value = -14808
$ws.Range("N28").Value = -35384
$ws.Range("H102").Value = 1060.3334
$ws.Range("I102").Value = 752.4
$ws.Range("K102").Value = 752.4
$ws.Range("M102").Value = 869.6

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 451.42856
$ws.Range("I16").Value = 376.66666
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 376.66666
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = -206.66666
$ws.Range("N16").Value = -1240
$ws.Range("H22").Value = 2265.1667
$ws.Range("I22").Value = 681.3333
$ws.Range("K22").Value = 681.3333
$ws.Range("M22").Value = -386.3333
$ws.Range("H27").Value = 2265.1667
$ws.Range("I27").Value = 681.3333
$ws.Range("K27").Value = 681.3333
$ws.Range("M27").Value = -574.3333
$ws.Range("H82").Value = 1275.2354
$ws.Range("I82").Value = 535.0909
$ws.Range("K82").Value = 535.0909
$ws.Range("M82").Value = -174.0909
$ws.Range("H85").Value = 1275.2354
$ws.Range("I85").Value = 535.0909
$ws.Range("K85").Value = 535.0909
$ws.Range("M85").Value = 712.9091
$ws.Range("H93").Value = 2995.4167
$ws.Range("I93").Value = 2811.125
$ws.Range("K93").Value = 2811.125
$ws.Range("M93").Value = -1563.125
$ws.Range("H103").Value = 7664
$ws.Range("J103").Value = 7664
$ws.Range("L103").Value = 7664
$ws.Range("N103").Value = -10008

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 14333
$ws.Range("J80").Value = 14333
$ws.Range("L80").Value = 14333
$ws.Range("N80").Value = -16329
$ws.Range("H83").Value = 14333
$ws.Range("J83").Value = 14333
$ws.Range("L83").Value = 42999
$ws.Range("N83").Value = -52983
$ws.Range("H126").Value = 7373.5
$ws.Range("J126").Value = 6492.5
$ws.Range("L126").Value = 19477.5
$ws.Range("N126").Value = -24417.5
$ws.Range("H132").Value = 2676.9443
$ws.Range("I132").Value = 2446.6667
$ws.Range("J132").Value = 3828.3333
$ws.Range("K132").Value = 7340.000100000001
$ws.Range("L132").Value = 11484.9999
$ws.Range("M132").Value = -4810.000100000001
$ws.Range("N132").Value = -16544.9999
$ws.Range("H136").Value = 11330.72
$ws.Range("I136").Value = 10076.647
$ws.Range("J136").Value = 13995.625
$ws.Range("K136").Value = 30229.941
$ws.Range("L136").Value = 41986.875
$ws.Range("M136").Value = -27679.941
$ws.Range("N136").Value = -47086.875
